$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new "2022" column (S) by extending the existing table one column
#     to the right, copying the formatting from the last year column (R). ---

# Row 2 (thin separator row above the header) - empty cell, same style as R2.
$ws.Range("R2").Copy($ws.Range("S2")) | Out-Null

# Row 3 (year header row) - new year value 2022, same style as R3.
$ws.Range("R3").Copy($ws.Range("S3")) | Out-Null
$ws.Range("S3").Value = 2022

# Row 4 (GVA share %) - revise the last three years' figures and add 2022.
$ws.Range("P4").Value = 13.7
$ws.Range("Q4").Value = 13.1
$ws.Range("R4").Value = 11.8
$ws.Range("R4").Copy($ws.Range("S4")) | Out-Null
$ws.Range("S4").Value = 13.6

# Row 5 (GVA per capita) - revise the last three years' figures and add 2022.
$ws.Range("P5").Value = 13.6
$ws.Range("Q5").Value = 12.5
$ws.Range("R5").Value = 13.5
$ws.Range("R5").Copy($ws.Range("S5")) | Out-Null
$ws.Range("S5").Value = 20

# Match the author's final selection (cell S2, the freshly added blank cell).
$ws.Range("S2").Select() | Out-Null
